# The sheet has a table "Tabla1" (A3:C9) with columns: Codi | Nom | Codi diputació.
# The edit removes the "Codi diputació" column, keeping its (longer) codes as the
# values of the "Codi" column, and collapses the table/worksheet back down to two
# columns (Codi | Nom).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Move the "Codi diputació" values into the "Codi" column.
$lo.ListColumns.Item(1).DataBodyRange.Value = $lo.ListColumns.Item(3).DataBodyRange.Value()

# Drop the now-redundant "Codi diputació" column from the table definition.
$lo.ListColumns.Item(3).Delete()

# Physically delete the emptied worksheet column so later columns shift left.
$ws.Range("C:C").Delete()
